$wb = $excel.ActiveWorkbook

# Relabel the light/heavy goods vehicle categories to van/lorry
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $text = $cell.Text
        if ($text -eq "light goods") {
            $cell.Value = "van"
        } elseif ($text -eq "heavy goods") {
            $cell.Value = "lorry"
        }
    }
}

# Updated simulation output values per sheet

$ws = $wb.Worksheets.Item("mean")
$ws.Range("B2").Value = 4.4151396369313165
$ws.Range("C2").Value = 3.7711676559002907
$ws.Range("D2").Value = 5.4833909095509314
$ws.Range("E2").Value = 2.6434965643222883
$ws.Range("F2").Value = 3.931372644906297
$ws.Range("B3").Value = 3.2944168478011604
$ws.Range("C3").Value = 2.4386713509470788
$ws.Range("D3").Value = 3.030701270140836
$ws.Range("E3").Value = 1.7035782217459488
$ws.Range("F3").Value = 2.6181634790175554
$ws.Range("B4").Value = 16.041206213557164
$ws.Range("C4").Value = 19.102099995259564
$ws.Range("D4").Value = 20.097756706112666
$ws.Range("E4").Value = 18.70197355391261
$ws.Range("F4").Value = 17.27223623696062
$ws.Range("B5").Value = 10.057335501557677
$ws.Range("C5").Value = 12.428766572404774
$ws.Range("D5").Value = 9.868352661393699
$ws.Range("E5").Value = 4.057077077101366
$ws.Range("F5").Value = 8.172104047970361
$ws.Range("B6").Value = 22.92401215404581
$ws.Range("C6").Value = 25.908772927729466
$ws.Range("D6").Value = 19.11737602999527
$ws.Range("E6").Value = 13.532467741818603
$ws.Range("F6").Value = 19.48438223048522
$ws.Range("B7").Value = 5.108548652920153
$ws.Range("C7").Value = 1.710876233771123
$ws.Range("D7").Value = 1.4284341858201288
$ws.Range("E7").Value = 0.7350405790411106
$ws.Range("F7").Value = 1.2433651943739956

$ws = $wb.Worksheets.Item("median")
$ws.Range("B2").Value = 4.418551932322057
$ws.Range("C2").Value = 3.76785896876474
$ws.Range("D2").Value = 5.480867972620296
$ws.Range("E2").Value = 2.6422637369238355
$ws.Range("F2").Value = 3.9306719774828194
$ws.Range("B3").Value = 3.293631446211991
$ws.Range("C3").Value = 2.436673938523467
$ws.Range("D3").Value = 3.027615205970023
$ws.Range("E3").Value = 1.695611644566232
$ws.Range("F3").Value = 2.6190257589582533
$ws.Range("B4").Value = 16.043201101125227
$ws.Range("C4").Value = 19.084521280358956
$ws.Range("D4").Value = 20.05201832024916
$ws.Range("E4").Value = 18.65449772825528
$ws.Range("F4").Value = 17.27467465279599
$ws.Range("B5").Value = 10.008192138881967
$ws.Range("C5").Value = 12.391306948855176
$ws.Range("D5").Value = 9.737401844433142
$ws.Range("E5").Value = 4.014620976797044
$ws.Range("F5").Value = 8.145151009500353
$ws.Range("B6").Value = 22.82602198908156
$ws.Range("C6").Value = 25.918629633632357
$ws.Range("D6").Value = 18.938390694368273
$ws.Range("E6").Value = 13.462031425076564
$ws.Range("F6").Value = 19.447033824724627
$ws.Range("B7").Value = 4.71371174920569
$ws.Range("C7").Value = 1.644084861902837
$ws.Range("D7").Value = 1.383934485398297
$ws.Range("E7").Value = 0.7136738556469473
$ws.Range("F7").Value = 1.2268040804629687

$ws = $wb.Worksheets.Item("lower 5")
$ws.Range("B2").Value = 4.266719367012875
$ws.Range("C2").Value = 3.6191421353928646
$ws.Range("D2").Value = 5.245471786810901
$ws.Range("E2").Value = 2.5297845787697533
$ws.Range("F2").Value = 3.8551490435672333
$ws.Range("B3").Value = 3.066102502571442
$ws.Range("C3").Value = 2.1840118993747386
$ws.Range("D3").Value = 2.742002559677885
$ws.Range("E3").Value = 1.5232659332529144
$ws.Range("F3").Value = 2.497494700024752
$ws.Range("B4").Value = 15.372300865637502
$ws.Range("C4").Value = 17.635442209929543
$ws.Range("D4").Value = 18.060052772558475
$ws.Range("E4").Value = 16.841135874360685
$ws.Range("F4").Value = 16.67190943312802
$ws.Range("B5").Value = 8.43813475584457
$ws.Range("C5").Value = 10.473558876174668
$ws.Range("D5").Value = 8.011467311445735
$ws.Range("E5").Value = 3.300456455726389
$ws.Range("F5").Value = 7.433569627310622
$ws.Range("B6").Value = 20.049192766056297
$ws.Range("C6").Value = 23.454676545599384
$ws.Range("D6").Value = 15.923798033310396
$ws.Range("E6").Value = 12.046243506455406
$ws.Range("F6").Value = 18.347179656565253
$ws.Range("B7").Value = 2.5809592257216276
$ws.Range("C7").Value = 1.0293818597991533
$ws.Range("D7").Value = 0.8409661051329211
$ws.Range("E7").Value = 0.466845365680859
$ws.Range("F7").Value = 0.9524821835445497

$ws = $wb.Worksheets.Item("upper 95")
$ws.Range("B2").Value = 4.555772214326913
$ws.Range("C2").Value = 3.937156551462659
$ws.Range("D2").Value = 5.716783119371994
$ws.Range("E2").Value = 2.7613289939541223
$ws.Range("F2").Value = 4.007388685485368
$ws.Range("B3").Value = 3.5438870682425923
$ws.Range("C3").Value = 2.719872015415003
$ws.Range("D3").Value = 3.327953066149144
$ws.Range("E3").Value = 1.8978098847771496
$ws.Range("F3").Value = 2.7386207473545823
$ws.Range("B4").Value = 16.822480440592397
$ws.Range("C4").Value = 20.571112613194558
$ws.Range("D4").Value = 22.36174451515825
$ws.Range("E4").Value = 20.68881896959167
$ws.Range("F4").Value = 17.836081448167665
$ws.Range("B5").Value = 11.85163195251912
$ws.Range("C5").Value = 14.512923056338273
$ws.Range("D5").Value = 12.026163322696128
$ws.Range("E5").Value = 4.936649176817264
$ws.Range("F5").Value = 8.971049829035898
$ws.Range("B6").Value = 26.018705149590925
$ws.Range("C6").Value = 28.319531841253394
$ws.Range("D6").Value = 22.670027694754136
$ws.Range("E6").Value = 15.135221859071812
$ws.Range("F6").Value = 20.740149896034264
$ws.Range("B7").Value = 9.028242037502663
$ws.Range("C7").Value = 2.6526229894777473
$ws.Range("D7").Value = 2.178907047210746
$ws.Range("E7").Value = 1.0719872902324334
$ws.Range("F7").Value = 1.6002729250966354
